$d = $word.ActiveDocument

$d.Content.Find.Execute("149×3=447", $true, $false, $false, $false, $false, $true, 1, $false, "489×8=3912", 2) | Out-Null
$d.Content.Find.Execute("641×3=1923", $true, $false, $false, $false, $false, $true, 1, $false, "564×6=3384", 2) | Out-Null
$d.Content.Find.Execute("282×3=846", $true, $false, $false, $false, $false, $true, 1, $false, "832×5=4160", 2) | Out-Null
$d.Content.Find.Execute("671×3=2013", $true, $false, $false, $false, $false, $true, 1, $false, "534×2=1068", 2) | Out-Null
$d.Content.Find.Execute("436×2=872", $true, $false, $false, $false, $false, $true, 1, $false, "978×4=3912", 2) | Out-Null
$d.Content.Find.Execute("863×2=1726", $true, $false, $false, $false, $false, $true, 1, $false, "889×9=8001", 2) | Out-Null
$d.Content.Find.Execute("168×2=336", $true, $false, $false, $false, $false, $true, 1, $false, "165×4=660", 2) | Out-Null
$d.Content.Find.Execute("686×7=4802", $true, $false, $false, $false, $false, $true, 1, $false, "254×5=1270", 2) | Out-Null
$d.Content.Find.Execute("945×9=8505", $true, $false, $false, $false, $false, $true, 1, $false, "775×3=2325", 2) | Out-Null
$d.Content.Find.Execute("693×8=5544", $true, $false, $false, $false, $false, $true, 1, $false, "575×8=4600", 2) | Out-Null
$d.Content.Find.Execute("478×8=3824", $true, $false, $false, $false, $false, $true, 1, $false, "689×8=5512", 2) | Out-Null
$d.Content.Find.Execute("991×2=1982", $true, $false, $false, $false, $false, $true, 1, $false, "478×7=3346", 2) | Out-Null
$d.Content.Find.Execute("380×2=760", $true, $false, $false, $false, $false, $true, 1, $false, "888×6=5328", 2) | Out-Null
$d.Content.Find.Execute("979×9=8811", $true, $false, $false, $false, $false, $true, 1, $false, "418×4=1672", 2) | Out-Null
$d.Content.Find.Execute("297×8=2376", $true, $false, $false, $false, $false, $true, 1, $false, "181×4=724", 2) | Out-Null
$d.Content.Find.Execute("155×7=1085", $true, $false, $false, $false, $false, $true, 1, $false, "978×4=3912", 2) | Out-Null
$d.Content.Find.Execute("392×9=3528", $true, $false, $false, $false, $false, $true, 1, $false, "639×8=5112", 2) | Out-Null
$d.Content.Find.Execute("494×9=4446", $true, $false, $false, $false, $false, $true, 1, $false, "621×8=4968", 2) | Out-Null
$d.Content.Find.Execute("550×7=3850", $true, $false, $false, $false, $false, $true, 1, $false, "501×2=1002", 2) | Out-Null
$d.Content.Find.Execute("485×8=3880", $true, $false, $false, $false, $false, $true, 1, $false, "131×8=1048", 2) | Out-Null
$d.Content.Find.Execute("578×8=4624", $true, $false, $false, $false, $false, $true, 1, $false, "521×7=3647", 2) | Out-Null
$d.Content.Find.Execute("837×2=1674", $true, $false, $false, $false, $false, $true, 1, $false, "959×9=8631", 2) | Out-Null
$d.Content.Find.Execute("554×2=1108", $true, $false, $false, $false, $false, $true, 1, $false, "635×2=1270", 2) | Out-Null
$d.Content.Find.Execute("190×6=1140", $true, $false, $false, $false, $false, $true, 1, $false, "524×6=3144", 2) | Out-Null
$d.Content.Find.Execute("324×7=2268", $true, $false, $false, $false, $false, $true, 1, $false, "423×7=2961", 2) | Out-Null
